# Refresh the crypto price/volume snapshot (and restore the original
# row ordering for a handful of coins whose rank swapped) per the
# GitHub Actions scrape dated Wed May  1 19:10:03 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range("D2").Value = "'" + '58.974.56'
$ws.Range("E2").Value = '  -2.39%  '

# row 3
$ws.Range("D3").Value = "'" + '2.989.29'
$ws.Range("E3").Value = '  +0.34%  '

# row 4
$ws.Range("D4").Value = "'" + '1.00'
$ws.Range("E4").Value = '  +0.32%  '

# row 5
$ws.Range("D5").Value = "'" + '561.05'
$ws.Range("E5").Value = '  -2.54%  '

# row 6
$ws.Range("D6").Value = "'" + '133.37'
$ws.Range("E6").Value = '  +7.18%  '

# row 7
$ws.Range("D7").Value = "'" + '1.00'
$ws.Range("E7").Value = '  +0.23%  '

# row 8
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = "'" + '0.516'
$ws.Range("E8").Value = '  +3.62%  '

# row 9
$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D9").Value = "'" + '2.980.08'
$ws.Range("E9").Value = '  +0.18%  '

# row 10
$ws.Range("D10").Value = "'" + '0.131'
$ws.Range("E10").Value = '  -1.29%  '

# row 11
$ws.Range("D11").Value = "'" + '4.89'
$ws.Range("E11").Value = '  -4.30%  '

# row 12
$ws.Range("D12").Value = "'" + '0.453'
$ws.Range("E12").Value = '  +4.07%  '

# row 13
$ws.Range("D13").Value = "'" + '0.0000227'
$ws.Range("E13").Value = '  +2.19%  '

# row 14
$ws.Range("D14").Value = "'" + '33.35'
$ws.Range("E14").Value = '  +2.47%  '

# row 15
$ws.Range("E15").Value = '  +2.15%  '

# row 16
$ws.Range("D16").Value = "'" + '3.482.80'
$ws.Range("E16").Value = '  +0.31%  '

# row 17
$ws.Range("E17").Value = '  +11.69%  '

# row 18
$ws.Range("D18").Value = "'" + '2.990.79'
$ws.Range("E18").Value = '  +0.71%  '

# row 19
$ws.Range("D19").Value = "'" + '58.960.18'
$ws.Range("E19").Value = '  -2.30%  '

# row 20
$ws.Range("D20").Value = "'" + '425.60'
$ws.Range("E20").Value = '  -0.42%  '

# row 21
$ws.Range("D21").Value = "'" + '13.31'
$ws.Range("E21").Value = '  +2.17%  '

# row 22
$ws.Range("D22").Value = "'" + '0.691'
$ws.Range("E22").Value = '  +4.72%  '

# row 23
$ws.Range("D23").Value = "'" + '7.07'
$ws.Range("E23").Value = '  +0.43%  '

# row 24
$ws.Range("D24").Value = "'" + '13.21'
$ws.Range("E24").Value = '  +2.79%  '

# row 25
$ws.Range("D25").Value = "'" + '80.25'
$ws.Range("E25").Value = '  +2.10%  '

# row 26
$ws.Range("E26").Value = '  +0.13%  '

# row 27
$ws.Range("D27").Value = "'" + '1.00'
$ws.Range("E27").Value = '  +0.31%  '

# row 28
$ws.Range("D28").Value = "'" + '2.53'
$ws.Range("E28").Value = '  +0.30%  '

# row 29
$ws.Range("D29").Value = "'" + '7.71'
$ws.Range("E29").Value = '  +8.35%  '

# row 30
$ws.Range("D30").Value = "'" + '2.02'
$ws.Range("E30").Value = '  +8.13%  '

# row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = "'" + '0.107'
$ws.Range("E31").Value = '  +16.18%  '

# row 32
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = "'" + '6.27'
$ws.Range("E32").Value = '  +3.39%  '

# row 33
$ws.Range("D33").Value = "'" + '25.39'
$ws.Range("E33").Value = '  +0.79%  '

# row 34
$ws.Range("D34").Value = "'" + '2.16'
$ws.Range("E34").Value = '  -2.35%  '

# row 35
$ws.Range("B35").Value = 'Mantle'
$ws.Range("C35").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D35").Value = "'" + '0.958'
$ws.Range("E35").Value = '  +1.36%  '

# row 36
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").Value = "'" + '5.71'
$ws.Range("E36").Value = '  +2.79%  '

# row 37
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = "'" + '48.86'
$ws.Range("E37").Value = '  -1.21%  '

# row 38
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = "'" + '0.0₃0695'
$ws.Range("E38").Value = '  +6.68%  '

# row 39
$ws.Range("E39").Value = '  +8.29%  '

# row 40
$ws.Range("D40").Value = "'" + '2.62'
$ws.Range("E40").Value = '  +9.64%  '

# row 41
$ws.Range("D41").Value = "'" + '0.111'
$ws.Range("E41").Value = '  +1.48%  '

# row 42
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = "'" + '0.0354'
$ws.Range("E42").Value = '  -0.49%  '

# row 43
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = "'" + '385.46'
$ws.Range("E43").Value = '  +1.53%  '

# row 44
$ws.Range("D44").Value = "'" + '2.678.51'
$ws.Range("E44").Value = '  +1.44%  '

# row 45
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").Value = "'" + '0.999'
$ws.Range("E45").Value = '  +0.01%  '

# row 46
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").Value = "'" + '0.244'
$ws.Range("E46").Value = '  +4.23%  '

# row 47
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").Value = "'" + '121.38'
$ws.Range("E47").Value = '  +1.74%  '

# row 48
$ws.Range("B48").Value = 'Fetch.AI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D48").Value = "'" + '2.02'
$ws.Range("E48").Value = '  +3.42%  '

# row 49
$ws.Range("E49").Value = '  +3.57%  '

# row 50
$ws.Range("D50").Value = "'" + '23.87'
$ws.Range("E50").Value = '  +2.87%  '

# row 51
$ws.Range("D51").Value = "'" + '2.03'
$ws.Range("E51").Value = '  +2.72%  '
